# Baseline_summary.xlsx - "updated 100% data results"
#
# Adds the new 100%-RGBD-with-output-depth columns (E/F) to the
# srn_cars_100% sheet, marks the best ("novel") value per metric row in
# bold on every srn_cars_* / co3d_w_background sheet (bold+blue on the
# 100% sheet, plain bold elsewhere), adds MAX/MIN "best of row" helper
# formulas on the 50% and co3d_w_background sheets, and updates the
# selections / active sheet to match the author's final view state.

$wb = $excel.ActiveWorkbook

$xlCenter = -4108

# ---------------------------------------------------------------------
# srn_cars_50% ("RGBD_SI"/"RGBD_NN" best-of-row highlight + MAX/MIN cols)
# Doing the plain-bold sheets BEFORE the bold+blue 100% sheet below keeps
# the new font/style indices in the same creation order as the original
# edit (font 1 = bold black, font 2 = bold blue).
# ---------------------------------------------------------------------
$ws50 = $wb.Worksheets.Item("srn_cars_50%")

$ws50.Range("B19").Font.Bold = $true
$ws50.Range("D20").Font.Bold = $true
$ws50.Range("D21").Font.Bold = $true

$ws50.Range("G19").Formula = "=MAX(B19:F19)"
$ws50.Range("G20").Formula = "=MAX(B20:F20)"
$ws50.Range("G21").Formula = "=MIN(B21:F21)"

# ---------------------------------------------------------------------
# srn_cars_20% best-of-row highlight
# ---------------------------------------------------------------------
$ws20 = $wb.Worksheets.Item("srn_cars_20%")

$ws20.Range("D19").Font.Bold = $true
$ws20.Range("D20").Font.Bold = $true
$ws20.Range("B21").Font.Bold = $true

# ---------------------------------------------------------------------
# co3d_w_background best-of-row highlight + MAX/MIN helper column
# ---------------------------------------------------------------------
$wsCo3d = $wb.Worksheets.Item("co3d_w_background")

$wsCo3d.Range("B19").Font.Bold = $true
$wsCo3d.Range("B20").Font.Bold = $true
$wsCo3d.Range("B21").Font.Bold = $true

$wsCo3d.Range("I19").Formula = "=MAX(B19:C19)"
$wsCo3d.Range("I20").Formula = "=MAX(B20:C20)"
$wsCo3d.Range("I21").Formula = "=MIN(B21:C21)"

# ---------------------------------------------------------------------
# srn_cars_100% - new "RGBD_SI_OUT_DEPTH"/"RGBD_NN_OUT_DEPTH" columns
# (E/F) plus bold+blue best-of-row highlight
# ---------------------------------------------------------------------
$ws100 = $wb.Worksheets.Item("srn_cars_100%")

# All new E/F cells (rows 12-24) share the same centered style as the
# existing B/C/D columns - apply it to the whole block up front so every
# cell (including the blank spacer cells on rows 15/22) ends up styled.
$ws100.Range("E12:F24").HorizontalAlignment = $xlCenter

$ws100.Range("E12").Value = -1.9135
$ws100.Range("F12").Value = -1.9183300000000001

$ws100.Range("E13").Value = -1.84918
$ws100.Range("F13").Value = -1.85304

$ws100.Range("E14").Value = -0.68400000000000005
$ws100.Range("F14").Value = -0.68225000000000002

$ws100.Range("E16").Value = 24.737295589663699
$ws100.Range("F16").Value = 24.7535391118038

$ws100.Range("E17").Value = 0.93629159058698197
$ws100.Range("F17").Value = 0.93711998567662402

$ws100.Range("E18").Value = 0.15091960154346701
$ws100.Range("F18").Value = 0.14831589974081899

$ws100.Range("E19").Value = 18.8896101027761
$ws100.Range("F19").Value = 18.946477096588701

$ws100.Range("E20").Value = 0.82427795599265596
$ws100.Range("F20").Value = 0.82487530109726404

$ws100.Range("E21").Value = 0.26353628302848098
$ws100.Range("F21").Value = 0.26368653031844902

$ws100.Range("E23").Value = 60
$ws100.Range("F23").Value = 50

$ws100.Range("E24").Value = 36
$ws100.Range("F24").Value = 38

$ws100.Range("B19").Font.Bold = $true
$ws100.Range("B19").Font.Color = 12611584

$ws100.Range("D20").Font.Bold = $true
$ws100.Range("D20").Font.Color = 12611584

$ws100.Range("D21").Font.Bold = $true
$ws100.Range("D21").Font.Color = 12611584

# ---------------------------------------------------------------------
# View-state: selections on each sheet, and the final active tab/sheet.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$ws1.Range("I19").Select()

$ws50.Activate()
$ws50.Range("G22").Select()

$ws20.Activate()
$ws20.Range("E16").Select()

$wsCo3d.Activate()
$wsCo3d.Range("F14").Select()

$ws100.Activate()
$ws100.Range("I24").Select()
